$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarterly period headers (shift left by one quarter; oldest quarter
# 1399/06 drops off and the newest quarter 1401/12 is appended at the end)
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: report publish dates (same left-shift, newest date appended at the end).
# NOTE: plain "1401-10-28" (no trailing " (n)") looks like a literal date to Excel
# and would otherwise be silently converted into a date serial number, so for that
# one cell we force text formatting first and then restore the original cell format
# (via copy/paste-special-formats) so the style index is left untouched.
$ws.Range("D9").Value = "1400-11-04 (3)"
$ws.Range("E9").Value = "1401-04-04 (9)"
$ws.Range("F9").Value = "1401-05-05 (3)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-10-28 (3)"
$ws.Range("I9").Value = "1402-02-27 (9)"
$ws.Range("J9").Value = "1401-05-05 (2)"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("M9").Value = "1402-02-27 (2)"

# Rows 11-27 (financial figures): shift every quarterly column one to the left and
# append the new 1401/12 quarter value in column M (some older figures were also
# revised slightly due to the read_price algorithm change).
# Row 11
$ws.Range("D11").Value = 6886433
$ws.Range("E11").Value = 7197439
$ws.Range("F11").Value = 9287246
$ws.Range("G11").Value = 10388296
$ws.Range("H11").Value = 11093924
$ws.Range("I11").Value = 30345171
$ws.Range("J11").Value = 17892661
$ws.Range("K11").Value = 18712197
$ws.Range("L11").Value = 20167635
$ws.Range("M11").Value = 22509981

# Row 12
$ws.Range("D12").Value = -2231547
$ws.Range("E12").Value = -2997654
$ws.Range("F12").Value = -2578075
$ws.Range("G12").Value = -3462313
$ws.Range("H12").Value = -5870272
$ws.Range("I12").Value = -9544641
$ws.Range("J12").Value = -4812766
$ws.Range("K12").Value = -6017101
$ws.Range("L12").Value = -5431361
$ws.Range("M12").Value = -9212042

# Row 13
$ws.Range("D13").Value = 4654886
$ws.Range("E13").Value = 4199785
$ws.Range("F13").Value = 6709171
$ws.Range("G13").Value = 6925983
$ws.Range("H13").Value = 5223652
$ws.Range("I13").Value = 20800530
$ws.Range("J13").Value = 13079895
$ws.Range("K13").Value = 12695096
$ws.Range("L13").Value = 14736274
$ws.Range("M13").Value = 13297939

# Row 14
$ws.Range("D14").Value = -408851
$ws.Range("E14").Value = -484739
$ws.Range("F14").Value = -578081
$ws.Range("G14").Value = -543289
$ws.Range("H14").Value = -518742
$ws.Range("I14").Value = -619033
$ws.Range("J14").Value = -973166
$ws.Range("K14").Value = -766949
$ws.Range("L14").Value = -952232
$ws.Range("M14").Value = -1072824

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = -192196
$ws.Range("E16").Value = 259670
$ws.Range("F16").Value = -262799
$ws.Range("G16").Value = -519489
$ws.Range("H16").Value = 1149780
$ws.Range("I16").Value = 12064
$ws.Range("J16").Value = 1077920
$ws.Range("K16").Value = 189200
$ws.Range("L16").Value = 493664
$ws.Range("M16").Value = 1705544

# Row 17
$ws.Range("D17").Value = 4053839
$ws.Range("E17").Value = 3974716
$ws.Range("F17").Value = 5868291
$ws.Range("G17").Value = 5863205
$ws.Range("H17").Value = 5854690
$ws.Range("I17").Value = 20193561
$ws.Range("J17").Value = 13184649
$ws.Range("K17").Value = 12117347
$ws.Range("L17").Value = 14277706
$ws.Range("M17").Value = 13930659

# Row 18
$ws.Range("D18").Value = -49073
$ws.Range("E18").Value = -67934
$ws.Range("F18").Value = -70253
$ws.Range("G18").Value = -74176
$ws.Range("H18").Value = -65262
$ws.Range("I18").Value = -46558
$ws.Range("J18").Value = -52488
$ws.Range("K18").Value = -52487
$ws.Range("L18").Value = -39699
$ws.Range("M18").Value = -21945

# Row 19
$ws.Range("D19").Value = -976337
$ws.Range("E19").Value = -768134
$ws.Range("F19").Value = -544702
$ws.Range("G19").Value = 488850
$ws.Range("H19").Value = 1845173
$ws.Range("I19").Value = -2774621
$ws.Range("J19").Value = 79555
$ws.Range("K19").Value = -802514
$ws.Range("L19").Value = 675223
$ws.Range("M19").Value = -6714029

# Row 20
$ws.Range("D20").Value = 3028429
$ws.Range("E20").Value = 3138648
$ws.Range("F20").Value = 5253336
$ws.Range("G20").Value = 6277879
$ws.Range("H20").Value = 7634601
$ws.Range("I20").Value = 17372382
$ws.Range("J20").Value = 13211716
$ws.Range("K20").Value = 11262346
$ws.Range("L20").Value = 14913230
$ws.Range("M20").Value = 7194685

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 2743063
$ws.Range("H21").Value = -2743063
$ws.Range("I21").Value = 2743063
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("D22").Value = 3028429
$ws.Range("E22").Value = 3138648
$ws.Range("F22").Value = 5253336
$ws.Range("G22").Value = 9020942
$ws.Range("H22").Value = 4891538
$ws.Range("I22").Value = 20115445
$ws.Range("J22").Value = 13211716
$ws.Range("K22").Value = 11262346
$ws.Range("L22").Value = 14913230
$ws.Range("M22").Value = 7194685

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 3028429
$ws.Range("E24").Value = 3138648
$ws.Range("F24").Value = 5253336
$ws.Range("G24").Value = 9020942
$ws.Range("H24").Value = 4891538
$ws.Range("I24").Value = 20115445
$ws.Range("J24").Value = 13211716
$ws.Range("K24").Value = 11262346
$ws.Range("L24").Value = 14913230
$ws.Range("M24").Value = 7194685

# Row 25
$ws.Range("D25").Value = 858
$ws.Range("E25").Value = 889
$ws.Range("F25").Value = 1489
$ws.Range("G25").Value = 2556
$ws.Range("H25").Value = 1386
$ws.Range("I25").Value = 5700
$ws.Range("J25").Value = 3744
$ws.Range("K25").Value = 3191
$ws.Range("L25").Value = 4226
$ws.Range("M25").Value = 2039

# Row 26
$ws.Range("D26").Value = 3529200
$ws.Range("E26").Value = 3529200
$ws.Range("F26").Value = 3529200
$ws.Range("G26").Value = 3529200
$ws.Range("H26").Value = 3529200
$ws.Range("I26").Value = 3529200
$ws.Range("J26").Value = 3529200
$ws.Range("K26").Value = 3529200
$ws.Range("L26").Value = 3529200
$ws.Range("M26").Value = 3529200

# Row 27
$ws.Range("D27").Value = 858
$ws.Range("E27").Value = 889
$ws.Range("F27").Value = 1489
$ws.Range("G27").Value = 2556
$ws.Range("H27").Value = 1386
$ws.Range("I27").Value = 5700
$ws.Range("J27").Value = 3744
$ws.Range("K27").Value = 3191
$ws.Range("L27").Value = 4226
$ws.Range("M27").Value = 2039
